# Generate Report for Handoff
# Adds a new row (row 3) to the "Overview", "zh-cn" and "de-de" sheets for the
# newly handed-off file (GUID 99cb8a61-383a-48ad-ac26-c0077397b8ac), mirroring
# the existing row for c83b7686-6797-4458-bbbe-9ac7d779bc20 but with status
# "Ready for handoff" and the new handoff timestamps/targets.

$wb = $excel.ActiveWorkbook

$guidOld = "c83b7686-6797-4458-bbbe-9ac7d779bc20"
$guidNew = "99cb8a61-383a-48ad-ac26-c0077397b8ac"

$suffixMd = ""
for ($i = 0; $i -lt 149; $i++) { $suffixMd = $suffixMd + "o" }
$mdFile = "$guidNew$suffixMd.md"
$mdDisplay = "e2e\$mdFile"
$ghUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d83c88c3a5b83640c7444606d08b1276b122805/e2e/$mdFile"

$xlfSuffix = ""
for ($i = 0; $i -lt 42; $i++) { $xlfSuffix = $xlfSuffix + "o" }
$xlfHash = "e8e0e7128b7502127c8bcf00fefa72694a6a8adc"
$xlfZh = "$guidNew$xlfSuffix.$xlfHash.zh-cn.xlf"
$xlfDe = "$guidNew$xlfSuffix.$xlfHash.de-de.xlf"

$statusNew = "Ready for handoff"
$zhHandoffTime = "2016-08-15 14:27:31"
$deHandoffTime = "2016-08-15 14:27:35"

# ---------------------------------------------------------------------------
# Overview sheet (columns: File Name, Path And Name, Extension, Publish URL,
# zh-cn, de-de, Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Cells.Item(3, 1).Value = $mdFile
$wsOv.Cells.Item(3, 2).Value = $mdDisplay
$wsOv.Cells.Item(3, 3).Value = ".md"
$wsOv.Cells.Item(3, 4).Value = ""
$wsOv.Cells.Item(3, 5).Value = $statusNew
$wsOv.Cells.Item(3, 6).Value = $statusNew
$wsOv.Cells.Item(3, 7).Value = $deHandoffTime

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $ghUrl, "", "", $mdDisplay) | Out-Null

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (16 columns)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(3, 1).Value = $mdFile
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = $statusNew
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "False"
$wsZh.Cells.Item(3, 7).Value = $xlfZh
$wsZh.Cells.Item(3, 8).Value = $zhHandoffTime
$wsZh.Cells.Item(3, 9).Value = ""
$wsZh.Cells.Item(3, 10).Value = ""
$wsZh.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(3, 12).Value = ""
$wsZh.Cells.Item(3, 13).Value = "True"
$wsZh.Cells.Item(3, 14).Value = ""
$wsZh.Cells.Item(3, 15).Value = "False"
$wsZh.Cells.Item(3, 16).Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ghUrl, "", "", $mdFile) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (16 columns)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(3, 1).Value = $mdFile
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = $statusNew
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "False"
$wsDe.Cells.Item(3, 7).Value = $xlfDe
$wsDe.Cells.Item(3, 8).Value = $deHandoffTime
$wsDe.Cells.Item(3, 9).Value = ""
$wsDe.Cells.Item(3, 10).Value = ""
$wsDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(3, 12).Value = ""
$wsDe.Cells.Item(3, 13).Value = "True"
$wsDe.Cells.Item(3, 14).Value = ""
$wsDe.Cells.Item(3, 15).Value = "False"
$wsDe.Cells.Item(3, 16).Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ghUrl, "", "", $mdFile) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Column widths: the "Status" column (zh-cn/de-de col C) and the
# corresponding "zh-cn"/"de-de" columns on Overview (E/F) grow to fit the
# new, longer "Ready for handoff" text.
# ---------------------------------------------------------------------------
$wsOv.Columns.Item(5).AutoFit() | Out-Null
$wsOv.Columns.Item(6).AutoFit() | Out-Null
$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(3).AutoFit() | Out-Null

Write-Output "Report generated for handoff"
